$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$alokPath = "Alok/Results/Regression/Run_10-Apr-2020_07-01-01_PM/Excel Results"

# Row -> Results value (Passed/Failed). RunByName and RunInFolder are the
# same for every row that gets filled in.
$rows = @{
    12 = "Passed"
    13 = "Passed"
    14 = "Passed"
    15 = "Passed"
    16 = "Passed"
    18 = "Failed"
    19 = "Passed"
    20 = "Passed"
    21 = "Failed"
    22 = "Passed"
    23 = "Passed"
    25 = "Passed"
    26 = "Passed"
}

foreach ($r in $rows.Keys) {
    $ws.Range("B$r").Value = $rows[$r]
    $ws.Range("C$r").Value = "Alok"
    $ws.Range("D$r").Value = $alokPath
}
